# Added "Ser/Bld" to the "Intravascular- Any" grouper
#
# The LISTS sheet has a block of rows (32-42 originally) that all belong to
# the CHEM / SYSTEM / "Intravascular -any" grouper, listing valid SYSTEM
# values (Bld, BldA, BldC, BldMV, BldV, Ser, Plas, Ser/Plas, Ser/Pls/Bld,
# Ser/Plas.ultracentrifugate, BldP). A new value "Ser/Bld" needs to be
# inserted right after "Plas" (row 38) and before "Ser/Plas" (row 39),
# pushing every following row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 39 (shifts rows 39.. down to 40..,
# inheriting the formatting of the row above, which already has the
# correct style for the "Intravascular -any" group).
$ws.Rows(39).Insert()

# Populate the new row with the new grouper member.
$ws.Cells.Item(39, 1).Value = "CHEM"
$ws.Cells.Item(39, 2).Value = "SYSTEM"
$ws.Cells.Item(39, 3).Value = "Intravascular -any"
$ws.Cells.Item(39, 4).Value = "Ser/Bld"

# Reflect the author's final cursor position in the saved view state.
$ws.Range("D40").Select()
